$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 333, shifting rows 333:407 down to 334:408
$ws.Rows("333:333").Insert()

# Populate the newly inserted row 333 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,Q,R mirror the rest of this dataset (same
# market/product), only D,J,K,L,M,O,P are new values per the edit.
$ws.Range("A333").Value = 6
$ws.Range("B333").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C333").Value = 'Metropolitana'
$ws.Range("D333").Value = Get-Date -Year 2022 -Month 3 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("E333").Value = 13
$ws.Range("F333").Value = 100112043
$ws.Range("G333").Value = 'Pepino ensalada'
$ws.Range("H333").Value = 'Sin especificar'
$ws.Range("I333").Value = 'Primera'
$ws.Range("J333").Value = 500
$ws.Range("K333").Value = 16000
$ws.Range("L333").Value = 17000
$ws.Range("M333").Value = 16560
$ws.Range("N333").Value = '$/caja 60 unidades'
$ws.Range("O333").Value = 'Región de Arica y Parinacota'
$ws.Range("P333").Value = 276
$ws.Range("Q333").Value = 60
$ws.Range("R333").Value = 'Hortaliza'
